$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (liver.nii.gz): insert a new "02 Rigid MI" column before the
# old "02 B-Spline MI" column, which shifts to column E and is renamed
# "03 B-Spline MI". ---
$ws1.Columns.Item(4).Insert()
$ws1.Columns.Item(4).ColumnWidth = 11.14

$ws1.Range("D1").Value2 = "02 Rigid MI"
$ws1.Range("E1").Value2 = "03 B-Spline MI"

$ws1.Range("D2").Value2 = 0.84023
$ws1.Range("E2").Value2 = 0.88995
$ws1.Range("D3").Value2 = 0.81915
$ws1.Range("E3").Value2 = 0.88184
$ws1.Range("D4").Value2 = 0.86709
$ws1.Range("E4").Value2 = 0.91852
$ws1.Range("D5").Value2 = 0.91931
$ws1.Range("E5").Value2 = 0.93473
$ws1.Range("D6").Value2 = 0.88826
$ws1.Range("E6").Value2 = 0.93888
$ws1.Range("D7").Value2 = 0.80561
$ws1.Range("E7").Value2 = 0.829
$ws1.Range("D8").Value2 = 0.79748
$ws1.Range("E8").Value2 = 0.89346
$ws1.Range("D9").Value2 = 0.72985
$ws1.Range("E9").Value2 = 0.79707
$ws1.Range("D10").Value2 = 0.85243
$ws1.Range("E10").Value2 = 0.88187
$ws1.Range("D11").Value2 = 0.80288
$ws1.Range("E11").Value2 = 0.88471
$ws1.Range("D12").Value2 = 0.72985
$ws1.Range("E12").Value2 = 0.79707
$ws1.Range("D13").Value2 = 0.91931
$ws1.Range("E13").Value2 = 0.93888
$ws1.Range("D14").Value2 = 0.8309541666666668
$ws1.Range("E14").Value2 = 0.882165
$ws1.Range("D15").Value2 = 0.8309541666666668
$ws1.Range("E15").Value2 = 0.88471

# D14 is a brand-new "Mean" cell next to the already-highlighted C14/E14 -
# copy the highlight formatting so the whole mean row stays consistent.
$ws1.Range("C14").Copy()
$ws1.Range("D14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Sheet2 (tumor.nii.gz): same column insertion/rename. ---
$ws2.Columns.Item(4).Insert()
$ws2.Columns.Item(4).ColumnWidth = 11.14

$ws2.Range("D1").Value2 = "02 Rigid MI"
$ws2.Range("E1").Value2 = "03 B-Spline MI"

$ws2.Range("D2").Value2 = 0.78174
$ws2.Range("E2").Value2 = 0.87101
$ws2.Range("D3").Value2 = 0.85332
$ws2.Range("E3").Value2 = 0.89632
$ws2.Range("D4").Value2 = 0.72139
$ws2.Range("E4").Value2 = 0.72887
$ws2.Range("D5").Value2 = 0.90041
$ws2.Range("E5").Value2 = 0.93697
$ws2.Range("D6").Value2 = 0.86004
$ws2.Range("E6").Value2 = 0.93523
$ws2.Range("D7").Value2 = 0.09468
$ws2.Range("E7").Value2 = 0.02221
$ws2.Range("D8").Value2 = 0.5277
$ws2.Range("E8").Value2 = 0.57258
$ws2.Range("D9").Value2 = 0.62212
$ws2.Range("E9").Value2 = 0.60126
$ws2.Range("D10").Value2 = 0.79953
$ws2.Range("E10").Value2 = 0.8211
$ws2.Range("D11").Value2 = 0.64121
$ws2.Range("E11").Value2 = 0.63976
$ws2.Range("D12").Value2 = 0.09468
$ws2.Range("E12").Value2 = 0.02221
$ws2.Range("D13").Value2 = 0.90041
$ws2.Range("E13").Value2 = 0.93697
$ws2.Range("D14").Value2 = 0.6497691666666666
$ws2.Range("E14").Value2 = 0.6653741666666666
$ws2.Range("D15").Value2 = 0.72139
$ws2.Range("E15").Value2 = 0.72887

$ws2.Range("C14").Copy()
$ws2.Range("D14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Selection / active-sheet bookkeeping: the edit leaves the cursor on
# D20 on both sheets, with "liver.nii.gz" (sheet 1) as the active tab
# instead of "tumor.nii.gz" (sheet 2). ---
$ws2.Range("D20").Select()
$ws1.Activate()
$ws1.Range("D20").Select()
